$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 240.63637
$ws.Range("I9").Value = 176.57143
$ws.Range("K9").Value = 176.57143
$ws.Range("M9").Value = -7.571429999999992
$ws.Range("H40").Value = 2088.8
$ws.Range("J40").Value = 2192.5
$ws.Range("L40").Value = 2192.5
$ws.Range("N40").Value = -2542.5
$ws.Range("H88").Value = 12423.846
$ws.Range("I88").Value = 2441.2
$ws.Range("J88").Value = 18663
$ws.Range("K88").Value = 2441.2
$ws.Range("L88").Value = 18663
$ws.Range("M88").Value = -2035.2
$ws.Range("N88").Value = -19475
$ws.Range("H91").Value = 12423.846
$ws.Range("I91").Value = 2441.2
$ws.Range("J91").Value = 18663
$ws.Range("K91").Value = 2441.2
$ws.Range("L91").Value = 18663
$ws.Range("M91").Value = -1037.2
$ws.Range("N91").Value = -21471
$ws.Range("H98").Value = 1384.0385
$ws.Range("I98").Value = 1432.5238
$ws.Range("J98").Value = 1180.4
$ws.Range("K98").Value = 1432.5238
$ws.Range("L98").Value = 1180.4
$ws.Range("M98").Value = 65.47620000000006
$ws.Range("N98").Value = -4176.4
$ws.Range("H122").Value = 1384.0385
$ws.Range("I122").Value = 1432.5238
$ws.Range("J122").Value = 1180.4
$ws.Range("K122").Value = 4297.5714
$ws.Range("L122").Value = 3541.2
$ws.Range("M122").Value = -1847.5714
$ws.Range("N122").Value = -8441.200000000001
$ws.Range("H129").Value = 1761.1666
$ws.Range("I129").Value = 1466.7273
$ws.Range("J129").Value = 5000
$ws.Range("K129").Value = 4400.1819
$ws.Range("L129").Value = 15000
$ws.Range("M129").Value = 599.8181000000004
$ws.Range("N129").Value = -25000

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4168126
$ws.Range("I61").Value = 4763322.5
$ws.Range("K61").Value = 4763322.5
$ws.Range("M61").Value = -4763110.5
$ws.Range("H136").Value = 4168126
$ws.Range("I136").Value = 4763322.5
$ws.Range("K136").Value = 14289967.5
$ws.Range("M136").Value = -14287417.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 496.33334
$ws.Range("I64").Value = 490
$ws.Range("J64").Value = 499.5
$ws.Range("K64").Value = 490
$ws.Range("L64").Value = 499.5
$ws.Range("M64").Value = -265
$ws.Range("N64").Value = -949.5
$ws.Range("H67").Value = 496.33334
$ws.Range("I67").Value = 490
$ws.Range("J67").Value = 499.5
$ws.Range("K67").Value = 490
$ws.Range("L67").Value = 499.5
$ws.Range("M67").Value = 290
$ws.Range("N67").Value = -2059.5
$ws.Range("H86").Value = 2029.25
$ws.Range("I86").Value = 2329.2222
$ws.Range("J86").Value = 1129.3334
$ws.Range("K86").Value = 2329.2222
$ws.Range("L86").Value = 1129.3334
$ws.Range("M86").Value = -1206.2222
$ws.Range("N86").Value = -3375.3334
$ws.Range("H89").Value = 2029.25
$ws.Range("I89").Value = 2329.2222
$ws.Range("J89").Value = 1129.3334
$ws.Range("K89").Value = 11646.111
$ws.Range("L89").Value = 5646.666999999999
$ws.Range("M89").Value = -6030.111000000001
$ws.Range("N89").Value = -16878.667
$ws.Range("H94").Value = 903.75
$ws.Range("I94").Value = 776.9231
$ws.Range("J94").Value = 1453.3334
$ws.Range("K94").Value = 776.9231
$ws.Range("L94").Value = 1453.3334
$ws.Range("M94").Value = -325.9231
$ws.Range("N94").Value = -2355.3334
$ws.Range("H134").Value = 558554.4399999999
$ws.Range("I134").Value = 530642.7
$ws.Range("J134").Value = 767892.5
$ws.Range("K134").Value = 1591928.1
$ws.Range("L134").Value = 2303677.5
$ws.Range("M134").Value = -1589393.1
$ws.Range("N134").Value = -2308747.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 4160.8
$ws.Range("I19").Value = 6701.6665
$ws.Range("J19").Value = 349.5
$ws.Range("K19").Value = 6701.6665
$ws.Range("L19").Value = 349.5
$ws.Range("M19").Value = -6531.6665
$ws.Range("N19").Value = -689.5
$ws.Range("H24").Value = 4160.8
$ws.Range("I24").Value = 6701.6665
$ws.Range("J24").Value = 349.5
$ws.Range("K24").Value = 6701.6665
$ws.Range("L24").Value = 349.5
$ws.Range("M24").Value = -6531.6665
$ws.Range("N24").Value = -689.5
$ws.Range("H62").Value = 3096.8
$ws.Range("I62").Value = 2667
$ws.Range("K62").Value = 2667
$ws.Range("M62").Value = -2043
$ws.Range("H65").Value = 3096.8
$ws.Range("I65").Value = 2667
$ws.Range("K65").Value = 13335
$ws.Range("M65").Value = -10215
$ws.Range("H99").Value = 3550.4707
$ws.Range("I99").Value = 3093.2222
$ws.Range("J99").Value = 4064.875
$ws.Range("K99").Value = 3093.2222
$ws.Range("L99").Value = 4064.875
$ws.Range("M99").Value = -1595.2222
$ws.Range("N99").Value = -7060.875
$ws.Range("H126").Value = 3550.4707
$ws.Range("I126").Value = 3093.2222
$ws.Range("J126").Value = 4064.875
$ws.Range("K126").Value = 9279.6666
$ws.Range("L126").Value = 12194.625
$ws.Range("M126").Value = -6809.6666
$ws.Range("N126").Value = -17134.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 79616.30499999999
$ws.Range("I139").Value = 101701.2
$ws.Range("K139").Value = 305103.6
$ws.Range("M139").Value = -299963.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 6402.4546
$ws.Range("J43").Value = 18521
$ws.Range("L43").Value = 18521
$ws.Range("N43").Value = -18823
$ws.Range("H70").Value = 5196.3335
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 5196.3335
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H126").Value = 619547.7
$ws.Range("I126").Value = 982473.4
$ws.Range("J126").Value = 2574
$ws.Range("K126").Value = 2947420.2
$ws.Range("L126").Value = 7722
$ws.Range("M126").Value = -2944950.2
$ws.Range("N126").Value = -12662

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4031
$ws.Range("I7").Value = 3818.182
$ws.Range("J7").Value = 4499.2
$ws.Range("K7").Value = 3818.182
$ws.Range("L7").Value = 4499.2
$ws.Range("M7").Value = -3706.182
$ws.Range("N7").Value = -4723.2
$ws.Range("H22").Value = 1198.7646
$ws.Range("I22").Value = 722.8
$ws.Range("K22").Value = 722.8
$ws.Range("M22").Value = -427.8
$ws.Range("H27").Value = 1198.7646
$ws.Range("I27").Value = 722.8
$ws.Range("K27").Value = 722.8
$ws.Range("M27").Value = -615.8
$ws.Range("H40").Value = 3839.8
$ws.Range("I40").Value = 3733.1667
$ws.Range("K40").Value = 3733.1667
$ws.Range("M40").Value = -3597.1667
$ws.Range("H46").Value = 7629.727
$ws.Range("I46").Value = 14705.8
$ws.Range("K46").Value = 14705.8
$ws.Range("M46").Value = -14517.8
$ws.Range("H122").Value = 5215.1177
$ws.Range("J122").Value = 6998.5557
$ws.Range("L122").Value = 20995.6671
$ws.Range("N122").Value = -25895.6671
$ws.Range("H126").Value = 4031
$ws.Range("I126").Value = 3818.182
$ws.Range("J126").Value = 4499.2
$ws.Range("K126").Value = 11454.546
$ws.Range("L126").Value = 13497.6
$ws.Range("M126").Value = -8984.545999999998
$ws.Range("N126").Value = -18437.6
$ws.Range("H132").Value = 970913.2
$ws.Range("I132").Value = 1292487.9
$ws.Range("J132").Value = 6189
$ws.Range("K132").Value = 3877463.7
$ws.Range("L132").Value = 18567
$ws.Range("M132").Value = -3874933.7
$ws.Range("N132").Value = -23627
$ws.Range("H136").Value = 3574.9092
$ws.Range("I136").Value = 2843.1333
$ws.Range("J136").Value = 5143
$ws.Range("K136").Value = 8529.3999
$ws.Range("L136").Value = 15429
$ws.Range("M136").Value = -5979.3999
$ws.Range("N136").Value = -20529

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3011.36
$ws.Range("I122").Value = 2472
$ws.Range("K122").Value = 7416
$ws.Range("M122").Value = -4966
$ws.Range("H126").Value = 4085.9614
$ws.Range("I126").Value = 3914.5652
$ws.Range("K126").Value = 11743.6956
$ws.Range("M126").Value = -9273.695599999999
$ws.Range("H136").Value = 16699.562
$ws.Range("I136").Value = 17212.7
$ws.Range("J136").Value = 9002.5
$ws.Range("K136").Value = 51638.10000000001
$ws.Range("L136").Value = 27007.5
$ws.Range("M136").Value = -49088.10000000001
$ws.Range("N136").Value = -32107.5
